$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column for rows 2-5 from 45174 to 45175
$ws.Range("C2:C5").Value = 45175
